$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (overwrite existing data) ---
$ws.Range("A3").Value = "313.213.213-2"
$ws.Range("B3").Value = "sdasda"
$ws.Range("C3").Value = "123@sdasd"
$ws.Range("D3").Value = "(13) 2132-1"
$ws.Range("E3").Value = "(32) 1321-3213"

# --- Row 4 ---
$ws.Range("A4").Value = 321.3
$ws.Range("B4").Value = "Rodrigo Bruno de Souza"
$ws.Range("C4").Value = "1@sdf.coms"
$ws.Range("D4").Value = "(12) 3132-13213"
$ws.Range("E4").Value = "(31) 321"

# --- Row 5 ---
$ws.Range("A5").Value = "321.321.321-32"
$ws.Range("B5").Value = "Rodrigo Bruno de Souza"
$ws.Range("C5").Value = "1321321@31321.com"
$ws.Range("D5").Value = "(31) 3213-21321"
$ws.Range("E5").Value = "(31) 3333-3333"

# --- Row 6 ---
$ws.Range("A6").Value = 3132
$ws.Range("B6").Value = "asd"
$ws.Range("C6").Value = "r@r.com"
$ws.Range("D6").Value = "(32) 1321-32132"
$ws.Range("E6").Value = "(32) 1"

# --- Row 7 ---
$ws.Range("A7").Value = "321.321.321"
$ws.Range("B7").Value = "sdfsdf"
$ws.Range("C7").Value = "1321321@asd"
$ws.Range("D7").Value = "(32) 1321-23132"
$ws.Range("E7").Value = "(32) 132"

# --- Row 8 ---
$ws.Range("A8").Value = "399.065.215-84"
$ws.Range("B8").Value = "Rodrigo Bruno de Souza"
$ws.Range("C8").Value = "ro_web1@fsa.com"
$ws.Range("D8").Value = "(12) 1321-32132"
$ws.Range("E8").Value = "(11) 9806-5588"

# --- Row 9 ---
$ws.Range("A9").Value = 33333333333333
$ws.Range("B9").Value = "asd"
$ws.Range("C9").Value = "a@sd.com"
$ws.Range("D9").Value = "(32) 1321-32132"
$ws.Range("E9").Value = "(32) 1321"

# --- Row 10 ---
$ws.Range("A10").Value = "313.212.313-2"
$ws.Range("B10").Value = "23a1d32asd1"
$ws.Range("C10").Value = "asd@asdas.com"
$ws.Range("D10").Value = "(32) 1321-321"
$ws.Range("E10").Value = "(32) 1321-3213"

# --- Row 11 ---
$ws.Range("A11").Value = "321.321.321-31"
$ws.Range("B11").Value = "32as1d3s2a1"
$ws.Range("C11").Value = "321321321321313@adsadas.com"
$ws.Range("D11").Value = "(32) 1321"
$ws.Range("E11").Value = "(31) 3213-21"

# --- Row 12 ---
$ws.Range("A12").Value = "313.213.213-21"
$ws.Range("B12").Value = "a23ds1s3a21"
$ws.Range("C12").Value = "313131@asddasd.com"
$ws.Range("D12").Value = "(32) 1"
$ws.Range("E12").Value = "(13) 21321"

# --- Row 13 (new) ---
$ws.Range("A13").Value = "313.213.213"
$ws.Range("B13").Value = 21321321
$ws.Range("C13").Value = "teste@dasdasd.com"
$ws.Range("D13").Value = "(32) 1"
$ws.Range("E13").Value = "(32) 1"

# --- Row 14 (new) ---
$ws.Range("A14").Value = "333.333.333-33"
$ws.Range("B14").Value = "asdsad"
$ws.Range("C14").Value = "asdsa@adsad.com"
$ws.Range("D14").Value = "(32) 1"
$ws.Range("E14").Value = "(65) 321"

# --- Row 15 (new) ---
$ws.Range("A15").Value = "313.213.21"
$ws.Range("B15").Value = "sssdfsd"
$ws.Range("C15").Value = "teste@sdsf.com"
$ws.Range("D15").Value = "(31) 3213-2113"
$ws.Range("E15").Value = "(32) 1321"

# --- Row 16 (new, A/D/E left blank) ---
$ws.Range("B16").Value = "asd"
$ws.Range("C16").Value = "3213132@sadasd"

# --- Row 17 (new) ---
$ws.Range("A17").Value = "333.333.333"
$ws.Range("B17").Value = "asd"
$ws.Range("C17").Value = "3213132@sadasd"
$ws.Range("D17").Value = "(31) 3213-2"
$ws.Range("E17").Value = "(32) 1321-321"

# --- Column width adjustments (B, C, E); col width model applies a
# +5/7 padding and rounds to the nearest 1/7 character unit, so we
# back the requested ColumnWidth off by 5/7 to land as close as
# possible to the target stored widths. ---
$ws.Columns.Item(2).ColumnWidth = 21.789887640449443 - 0.7142857142857143
$ws.Columns.Item(3).ColumnWidth = 25.08988764044944 - 0.7142857142857143
$ws.Columns.Item(5).ColumnWidth = 15.18988764044944 - 0.7142857142857143
